# "Change to Linux env"
# - Bump the addListItem test value from "ListLinuxI" to "ListLinuYA"
# - Bump the createUser auto-numbered test id from 2735 to 2736
# - Leave the UI in the state the author saved it in: addListItem tab
#   selected/active, createUser's last selection parked on D25.

$wb = $excel.ActiveWorkbook

$listSheet = $wb.Worksheets.Item("addListItem")
$listSheet.Range("A2").Value = "ListLinuYA"

$userSheet = $wb.Worksheets.Item("createUser")
$userSheet.Range("A2").Value = 2736

# Park the selection on createUser before switching away from it, matching
# the saved selection state.
$userSheet.Activate() | Out-Null
$userSheet.Range("D25").Select() | Out-Null

# addListItem ends up as the active/selected sheet.
$listSheet.Activate() | Out-Null
$listSheet.Range("A2").Select() | Out-Null
